$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for changed rows
# Leading apostrophe forces these to stay text cells (matching the
# original inlineStr/text storage) instead of being auto-converted to
# numbers/percentages by Excel value-parsing heuristics.
$ws.Range("D2").Value = "'309.22"
$ws.Range("E2").Value = "'1.13%"
$ws.Range("D3").Value = "'38.90"
$ws.Range("E3").Value = "'9.09%"
$ws.Range("D4").Value = "'5.105"
$ws.Range("E4").Value = "'1.29%"
$ws.Range("D5").Value = "'0.08135"
$ws.Range("E5").Value = "'1.78%"
$ws.Range("D6").Value = "'1.976"
$ws.Range("E6").Value = "'3.49%"
$ws.Range("D7").Value = "'7.927"
$ws.Range("E7").Value = "'2.10%"
$ws.Range("D8").Value = "'0.9296"
$ws.Range("E8").Value = "'0.93%"
$ws.Range("D9").Value = "'0.1432"
$ws.Range("E9").Value = "'12.27%"
$ws.Range("D10").Value = "'0.1957"
$ws.Range("E10").Value = "'1.71%"
$ws.Range("D11").Value = "'0.09216"
$ws.Range("E11").Value = "'1.14%"
$ws.Range("D12").Value = "'0.03489"
$ws.Range("E12").Value = "'1.23%"
$ws.Range("D13").Value = "'0.09834"
$ws.Range("E13").Value = "'0.02%"
$ws.Range("D14").Value = "'0.001409"
$ws.Range("E14").Value = "'0.66%"
$ws.Range("D15").Value = "'0.005821"
$ws.Range("E15").Value = "'-7.35%"
$ws.Range("D16").Value = "'3.594"
$ws.Range("E16").Value = "'-3.36%"
$ws.Range("D17").Value = "'4.192"
$ws.Range("E17").Value = "'1.05%"
$ws.Range("D19").Value = "'0.3446"
$ws.Range("E19").Value = "'-0.01%"
$ws.Range("D20").Value = "'0.1323"
$ws.Range("E20").Value = "'-1.41%"
$ws.Range("D21").Value = "'4.825"
$ws.Range("E21").Value = "'-6.54%"
$ws.Range("D22").Value = "'0.2467"
$ws.Range("E22").Value = "'-5.13%"
$ws.Range("D23").Value = "'0.04454"
$ws.Range("E23").Value = "'0.41%"
$ws.Range("D24").Value = "'0.001240"
$ws.Range("E24").Value = "'0.51%"
$ws.Range("D25").Value = "'0.004853"
$ws.Range("E25").Value = "'4.76%"
$ws.Range("D27").Value = "'0.0001303"
$ws.Range("D39").Value = "'0.02116"
$ws.Range("E39").Value = "'8.71%"
$ws.Range("D40").Value = "'0.05117"
$ws.Range("E40").Value = "'-3.92%"
$ws.Range("D41").Value = "'0.007481"
$ws.Range("E41").Value = "'-1.84%"
$ws.Range("D42").Value = "'0.009990"
$ws.Range("E42").Value = "'-2.00%"
$ws.Range("D43").Value = "'0.1364"
$ws.Range("E43").Value = "'0.63%"
$ws.Range("D44").Value = "'0.002145"
$ws.Range("D45").Value = "'0.01015"
$ws.Range("E45").Value = "'2.53%"
$ws.Range("D46").Value = "'0.00006245"
$ws.Range("E46").Value = "'1.75%"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("D49").Value = "'0.001603"
$ws.Range("E49").Value = "'-3.36%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.21%"

# Update Hora (G) column from 7 to 8 for all data rows (2-51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Range("G$r").Value = "'8"
}
